# Lesson 53 Updating an existing Note
#
# The document ends with a lone paragraph that only contains a single
# space character. Turn that into:
#   1) a bold, red "Note" style line about PUT vs POST for updates
#   2) four more plain paragraphs describing the three ways to pass data
#      (body, query string, route params), including the proofing marks
#      Word itself would add for the grammar/spelling quirks in the text
#   3) a trailing empty paragraph at the very end of the document
#
# We replace the content of that last paragraph (but not its paragraph
# mark) with a small WordprocessingML package via Range.InsertXML, which
# lets us specify the exact run/paragraph/proofErr structure in one shot.

$d = $word.ActiveDocument

$target = $d.Paragraphs.Last
$targetRange = $target.Range

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:color w:val="EE0000"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:color w:val="EE0000"/>
              </w:rPr>
              <w:t>You can do updation by post but for updation put is used.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>There ar</w:t>
            </w:r>
            <w:r>
              <w:t>e</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> three ways to get the data </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>Either by body in post</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">And by query </w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>in ?</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
          </w:p>
          <w:p>
            <w:r>
              <w:t>And params in /:&lt;</w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>parmaeter</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>&gt;</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$targetRange.InsertXML($xml)
